$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 613.4167
$ws.Range("I41").Value = 184
$ws.Range("K41").Value = 184
$ws.Range("M41").Value = 256
$ws.Range("H55").Value = 1084.9166
$ws.Range("I55").Value = 1403.625
$ws.Range("J55").Value = 447.5
$ws.Range("K55").Value = 1403.625
$ws.Range("L55").Value = 447.5
$ws.Range("M55").Value = -1189.625
$ws.Range("N55").Value = -875.5
$ws.Range("H74").Value = 10387325
$ws.Range("I74").Value = 10387325
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 10387325
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -10386389
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 10387325
$ws.Range("I77").Value = 10387325
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 51936625
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -51931945
$ws.Range("N77").ClearContents()
$ws.Range("H98").Value = 3157
$ws.Range("I98").Value = 3183.1667
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 3183.1667
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = -1685.1667
$ws.Range("N98").Value = -5996
$ws.Range("H116").Value = 5022
$ws.Range("I116").Value = 5116.3335
$ws.Range("J116").Value = 4833.3335
$ws.Range("K116").Value = 5116.3335
$ws.Range("L116").Value = 4833.3335
$ws.Range("M116").Value = -1674.3335
$ws.Range("N116").Value = -11717.3335
$ws.Range("H122").Value = 3157
$ws.Range("I122").Value = 3183.1667
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 9549.500100000001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -7099.500100000001
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 1682.75
$ws.Range("I132").Value = 1801
$ws.Range("J132").Value = 1117.7778
$ws.Range("K132").Value = 5403
$ws.Range("L132").Value = 3353.3334
$ws.Range("M132").Value = -2873
$ws.Range("N132").Value = -8413.3334
$ws.Range("H137").Value = 821.5294
$ws.Range("I137").Value = 725.1667
$ws.Range("J137").Value = 874.0909
$ws.Range("K137").Value = 2175.5001
$ws.Range("L137").Value = 2622.2727
$ws.Range("M137").Value = 374.4998999999998
$ws.Range("N137").Value = -7722.2727
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4091
$ws.Range("I32").Value = 3449.2651
$ws.Range("J32").Value = 7224.1763
$ws.Range("K32").Value = 3449.2651
$ws.Range("L32").Value = 7224.1763
$ws.Range("M32").Value = -3162.2651
$ws.Range("N32").Value = -7798.1763
$ws.Range("H132").Value = 1170.5161
$ws.Range("I132").Value = 1077.3334
$ws.Range("J132").Value = 1366.2
$ws.Range("K132").Value = 3232.0002
$ws.Range("L132").Value = 4098.6
$ws.Range("M132").Value = -702.0001999999999
$ws.Range("N132").Value = -9158.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 41693170
$ws.Range("I20").Value = 67798.336
$ws.Range("K20").Value = 67798.336
$ws.Range("M20").Value = -67551.336
$ws.Range("H134").Value = 17774.629
$ws.Range("I134").Value = 1428.7843
$ws.Range("J134").Value = 93559.91
$ws.Range("K134").Value = 4286.3529
$ws.Range("L134").Value = 280679.73
$ws.Range("M134").Value = -1751.3529
$ws.Range("N134").Value = -285749.73
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2496.9272
$ws.Range("I58").Value = 618.3143
$ws.Range("K58").Value = 618.3143
$ws.Range("M58").Value = -415.3143
$ws.Range("H94").Value = 2706.9211
$ws.Range("I94").Value = 3484.8572
$ws.Range("J94").Value = 2531.258
$ws.Range("K94").Value = 3484.8572
$ws.Range("L94").Value = 2531.258
$ws.Range("M94").Value = -3033.8572
$ws.Range("N94").Value = -3433.258
$ws.Range("H132").Value = 1590.4429
$ws.Range("I132").Value = 876.5814
$ws.Range("J132").Value = 2727.3333
$ws.Range("K132").Value = 2629.7442
$ws.Range("L132").Value = 8181.999899999999
$ws.Range("M132").Value = -99.74420000000009
$ws.Range("N132").Value = -13241.9999
$ws.Range("H134").Value = 1559.4762
$ws.Range("I134").Value = 1490.4242
$ws.Range("J134").Value = 1812.6666
$ws.Range("K134").Value = 4471.2726
$ws.Range("L134").Value = 5437.9998
$ws.Range("M134").Value = -1936.2726
$ws.Range("N134").Value = -10507.9998
$ws.Range("H136").Value = 2496.9272
$ws.Range("I136").Value = 618.3143
$ws.Range("K136").Value = 1854.9429
$ws.Range("M136").Value = 695.0571
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 563.4666999999999
$ws.Range("I113").Value = 579.1667
$ws.Range("J113").Value = 545.5238000000001
$ws.Range("K113").Value = 1737.5001
$ws.Range("L113").Value = 1636.5714
$ws.Range("M113").Value = 432.4999
$ws.Range("N113").Value = -5976.571400000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4089.4119
$ws.Range("I70").Value = 3977.2974
$ws.Range("J70").Value = 4385.7144
$ws.Range("K70").Value = 3977.2974
$ws.Range("L70").Value = 4385.7144
$ws.Range("M70").Value = -3707.2974
$ws.Range("N70").Value = -4925.7144
$ws.Range("H73").Value = 4089.4119
$ws.Range("I73").Value = 3977.2974
$ws.Range("J73").Value = 4385.7144
$ws.Range("K73").Value = 3977.2974
$ws.Range("L73").Value = 4385.7144
$ws.Range("M73").Value = -3041.2974
$ws.Range("N73").Value = -6257.7144
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6947760.5
$ws.Range("I7").Value = 3659.818
$ws.Range("K7").Value = 3659.818
$ws.Range("M7").Value = -3547.818
$ws.Range("H16").Value = 5275.793
$ws.Range("I16").Value = 7665.1763
$ws.Range("J16").Value = 1890.8334
$ws.Range("K16").Value = 7665.1763
$ws.Range("L16").Value = 1890.8334
$ws.Range("M16").Value = -7495.1763
$ws.Range("N16").Value = -2230.8334
$ws.Range("H40").Value = 632795.7
$ws.Range("I40").Value = 778325.0600000001
$ws.Range("K40").Value = 778325.0600000001
$ws.Range("M40").Value = -778189.0600000001
$ws.Range("H126").Value = 6947760.5
$ws.Range("I126").Value = 3659.818
$ws.Range("K126").Value = 10979.454
$ws.Range("M126").Value = -8509.454000000002
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 699.7308
$ws.Range("I132").Value = 591.2917
$ws.Range("K132").Value = 1773.8751
$ws.Range("M132").Value = 756.1249
$ws.Range("H136").Value = 856.7143
$ws.Range("I136").Value = 995.7037
$ws.Range("J136").Value = 606.5333000000001
$ws.Range("K136").Value = 2987.1111
$ws.Range("L136").Value = 1819.5999
$ws.Range("M136").Value = -437.1111000000001
$ws.Range("N136").Value = -6919.5999
